$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: a new time-log entry (row 121) was added, which pushed the trailing
# "blank helper" row and the "Total Time:" summary rows from 121-123 down to
# 151-153, leaving a block of still-unused placeholder rows (122-150) in
# between.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# 1) Grab the number format that the last trailing row (the one right under
#    the totals, "#\" days\"") uses, and stash a copy of it at its new home
#    (row 153) before anything below gets overwritten.
$ws1.Range("E123").Copy()
$ws1.Range("E153").PasteSpecial(-4122)

# 2) Remove the old "Total Time:" trailer block (rows 122 and 123) outright;
#    it gets rebuilt two rows further down.
$ws1.Range("A122").Clear()
$ws1.Range("E122").Clear()
$ws1.Range("E123").Clear()

# 3) Fill in the new logged entry on row 121.
$ws1.Range("A121").Value = 41968
$ws1.Range("B121").Value = 0.59791666666666665
$ws1.Range("C121").Value = 0.68819444444444444
$ws1.Range("D121").Value = 0
$ws1.Range("E121").Formula = '=IF(AND(NOT(ISBLANK(B121)),NOT(ISBLANK(C121))), (C121-B121) * 24 - D121/60, "")'
$ws1.Range("F121").Value = "Coding"

# 4) Carry the date-column formatting down through the new row and every
#    still-empty placeholder row that follows it.
$ws1.Range("A120").Copy()
$ws1.Range("A121:A150").PasteSpecial(-4122)

# 5) The helper formula row (evaluates to "" until more rows are filled in)
#    now lives at row 151.
$ws1.Range("E151").Formula = '=IF(AND(NOT(ISBLANK(B151)),NOT(ISBLANK(C151))), (C151-B151) * 24 - D151/60, "")'

# 6) Rebuild the "Total Time:" summary two rows further down, covering the
#    wider data range.
$ws1.Range("A152").Value = "Total Time:"
$ws1.Range("E152").Formula = "=SUM(E2:E151)"

# ---------------------------------------------------------------------------
# Sheet2: the SUMIF helper table has to look across the new, wider Sheet1
# data range (2:151 instead of 2:121).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B2").Formula = '=SUMIF(Sheet1!$F$2:$F$151,A2,Sheet1!$E$2:$E$151)'
$ws2.Range("B3").Formula = '=SUMIF(Sheet1!$F$2:$F$151,A3,Sheet1!$E$2:$E$151)'
$ws2.Range("B4").Formula = '=SUMIF(Sheet1!$F$2:$F$151,A4,Sheet1!$E$2:$E$151)'
$ws2.Range("B5").Formula = '=SUMIF(Sheet1!$F$2:$F$151,A5,Sheet1!$E$2:$E$151)'

# ---------------------------------------------------------------------------
# View-state bookkeeping: window size shrank and the active selection /
# scroll position on Sheet1 moved down along with the newly added row.
# ---------------------------------------------------------------------------
$excel.Width = 20730
$excel.Height = 11760

$ws1.Application.ActiveWindow.ScrollRow = 109
$ws1.Range("A122").Select()
